$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = 45949
$ws.Range("B2").Value = 0.009329675380027286
$ws.Range("C2").Value = 1.610176344153507
$ws.Range("D2").Value = 0.02747964137393155
$ws.Range("E2").Value = 0.00008704284269668728
$ws.Range("F2").Value = 29
$ws.Range("G2").Value = 5481.429871186317
$ws.Range("H2").Value = 74.03667922851697
$ws.Range("I2").Value = 55.92120153249586
$ws.Range("J2").Value = 0.7554013296252071
